# Update the "想去人数" (want-to-go count) column F values on the
# "展览" and "全部类型" worksheets to reflect the freshly generated
# site output (commit: "Update gh-pages to output generated at 7921097").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    3  = 7705
    6  = 39
    9  = 5927
    13 = 1807
    14 = 1321
    15 = 283
    16 = 548
    17 = 123
    18 = 5535
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
